# Daily refresh of the cryptos price list (GitHub Actions job).
# Column D (Price) and E (Volume 1h) are stored as text in this sheet, so
# numeric-looking prices are written with a leading apostrophe to stop
# Excel from auto-converting them to numbers (which would e.g. truncate
# "20.20" -> 20.2). The cell style is then reset to "Normal" so no stray
# text-format / quote-prefix style sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.201.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'1.647.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D5").Value = "'218.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'0.509"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.256"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "'0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'20.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'1.879.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'1.639.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "'4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "'0.537"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "'67.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "'27.182.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'0.0₃0738"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'220.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "'4.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("D24").Value = "'9.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'148.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'7.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").Value = "'0.0505"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'3.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "'1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'1.272.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("D38").Value = "'0.541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "'0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +6.74%  "
$ws.Range("D43").Value = "'5.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "'1.790.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'62.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").Value = "'92.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "'1.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.0₆0107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +16.65%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0972"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "
